# Regenerate save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates the "K" column (column G) values for rows 2-28 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 0
    6  = 2
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 2
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 1
    21 = 2
    22 = 0
    23 = 0
    24 = 1
    25 = 1
    26 = 0
    27 = 2
    28 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
